$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos list (Price column D and Volume(1h) column E) with the
# latest scraped values. A leading "'" forces numeric-looking price strings
# (e.g. "685.60", "1.00") to stay stored as exact text instead of being
# auto-converted to numbers and losing trailing zeros.
$ws.Range("D2").Value = "69.433.39"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "3.689.78"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'685.60"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'160.42"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.494"
$ws.Range("E8").Value = "  +2.08%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "'7.09"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "4.314.48"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'32.52"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "3.687.28"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "69.433.41"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'15.86"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").Value = "'470.58"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "'10.01"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "'79.87"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "3.837.60"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "'10.99"
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").Value = "'9.27"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Value = "'6.59"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").Value = "3.663.66"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("D37").Value = "'8.17"
$ws.Range("E37").Value = "  -3.19%  "
$ws.Range("D38").Value = "'6.16"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.24"
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'0.0901"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'166.14"
$ws.Range("E44").Value = "  +3.89%  "
$ws.Range("D45").Value = "'47.48"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("E46").Value = "  +6.87%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").Value = "'0.000281"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.73"
$ws.Range("E48").Value = "  -4.52%  "
$ws.Range("D49").Value = "'1.30"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").Value = "'28.21"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").Value = "'7.80"
$ws.Range("E51").Value = "  -1.85%  "
